$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "tags" column (O) mirroring the formatting of the existing
# "security" column (N): header style on row 1, data style on rows 2-4.

# O1 header "tags" - copy format from N1 (header style)
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("O1").Value = "tags"

# O2 "@elan @test" - copy format from N2
$ws.Range("N2").Copy()
$ws.Range("O2").PasteSpecial(-4122)
$ws.Range("O2").Value = "@elan @test"

# O3 "@pet @test" - copy format from N3
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)
$ws.Range("O3").Value = "@pet @test"

# O4 "@pet @test" - copy format from N3 (row4 has no N cell in source; use N3 style, same as O3)
$ws.Range("N3").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = "@pet @test"

$excel.CutCopyMode = 0

Write-Output "Added tags column with values for EDI271, PetGet and PetPost test cases"
